$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.101.76"
$ws.Range("E2").Value = "  -0.50%  "

$ws.Range("D3").Value = "1.799.81"
$ws.Range("E3").Value = "  -0.28%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.31%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.90"
$ws.Range("E5").Value = "  -1.11%  "

$ws.Range("E6").Value = "  -0.33%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5089"
$ws.Range("E7").Value = "  -2.98%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3866"
$ws.Range("E8").Value = "  +1.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07733"
$ws.Range("E9").Value = "  -2.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.096"
$ws.Range("E10").Value = "  -0.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.82"
$ws.Range("E11").Value = "  -2.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.328"
$ws.Range("E12").Value = "  -0.06%  "

$ws.Range("E13").Value = "  -0.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.23"
$ws.Range("E14").Value = "  -2.06%  "

$ws.Range("D15").Value = "1.800.90"
$ws.Range("E15").Value = "  -0.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.259"
$ws.Range("E16").Value = "  -1.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.95"
$ws.Range("E17").Value = "  -0.74%  "

$ws.Range("E18").Value = "  -1.73%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06561"
$ws.Range("E19").Value = "  -0.58%  "

$ws.Range("E20").Value = "  -0.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.19"
$ws.Range("E21").Value = "  -1.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.939"
$ws.Range("E22").Value = "  -0.54%  "

$ws.Range("D23").Value = "28.117.72"
$ws.Range("E23").Value = "  -0.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.04"
$ws.Range("E24").Value = "  -0.72%  "

$ws.Range("E25").Value = "  +0.58%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.70"
$ws.Range("E26").Value = "  +1.81%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.416"
$ws.Range("E27").Value = "  -0.29%  "

$ws.Range("D28").Value = "2.005.35"
$ws.Range("E28").Value = "  -0.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.16"
$ws.Range("E29").Value = "  -1.64%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.04"
$ws.Range("E30").Value = "  +3.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1089"
$ws.Range("E31").Value = "  -1.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.043"
$ws.Range("E32").Value = "  -1.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.647"
$ws.Range("E33").Value = "  -0.59%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.522"
$ws.Range("E34").Value = "  -1.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07025"
$ws.Range("E35").Value = "  -3.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.032"
$ws.Range("E36").Value = "  +2.93%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02340"
$ws.Range("E37").Value = "  +1.03%  "

$ws.Range("E38").Value = "  -0.77%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.015"
$ws.Range("E39").Value = "  -0.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.45"
$ws.Range("E40").Value = "  -5.62%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6089"
$ws.Range("E41").Value = "  -1.89%  "

$ws.Range("E42").Value = "  -0.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.151"
$ws.Range("E43").Value = "  -1.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.15"
$ws.Range("E44").Value = "  -0.82%  "

$ws.Range("E45").Value = "  -6.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5891"
$ws.Range("E46").Value = "  -2.75%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.713"
$ws.Range("E47").Value = "  -1.36%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.99"
$ws.Range("E48").Value = "  -0.86%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.194"
$ws.Range("E49").Value = "  -1.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.899"
$ws.Range("E50").Value = "  -1.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06741"
$ws.Range("E51").Value = "  -1.20%  "
